$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 335345.66
$ws.Range("J17").Value = 335345.66
$ws.Range("L17").Value = 1006036.98
$ws.Range("N17").Value = -1006372.98
# Row 19
$ws.Range("H19").Value = 948.3333
$ws.Range("J19").Value = 948.3333
$ws.Range("L19").Value = 948.3333
$ws.Range("N19").Value = -1298.3333
# Row 38
$ws.Range("H38").Value = 1872.8182
$ws.Range("I38").Value = 432.66666
$ws.Range("J38").Value = 3601.0
$ws.Range("K38").Value = 1297.99998
$ws.Range("L38").Value = 10803.0
$ws.Range("M38").Value = -925.99998
$ws.Range("N38").Value = -11547.0
# Row 39
$ws.Range("H39").Value = 5159.0
$ws.Range("I39").Value = 68.25
$ws.Range("K39").Value = 204.75
$ws.Range("M39").Value = 91.25
# Row 42
$ws.Range("H42").Value = 3232.0
$ws.Range("I42").Value = 3718.6667
$ws.Range("K42").Value = 11156.0001
$ws.Range("M42").Value = -10926.0001
# Row 43
$ws.Range("H43").Value = 2396.0
$ws.Range("I43").Value = 0.0
$ws.Range("J43").Value = 2396.0
$ws.Range("K43").Value = 0.0
$ws.Range("L43").Value = 2396.0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2534.0
# Row 70
$ws.Range("H70").Value = 6454.05
$ws.Range("I70").Value = 2630.0625
$ws.Range("K70").Value = 7890.1875
$ws.Range("M70").Value = -7620.1875
# Row 73
$ws.Range("H73").Value = 6454.05
$ws.Range("I73").Value = 2630.0625
$ws.Range("K73").Value = 7890.1875
$ws.Range("M73").Value = -6954.1875
# Row 96
$ws.Range("H96").Value = 2812.2222
$ws.Range("I96").Value = 2383.1667
$ws.Range("K96").Value = 7149.500100000001
$ws.Range("M96").Value = -5776.500100000001
# Row 103
$ws.Range("H103").Value = 1842.1428
$ws.Range("I103").Value = 1485.8334
$ws.Range("K103").Value = 4457.5002
$ws.Range("M103").Value = -3871.5002
# Row 112
$ws.Range("H112").Value = 66045.625
$ws.Range("I112").Value = 2690.0
$ws.Range("J112").Value = 70269.336
$ws.Range("K112").Value = 8070.0
$ws.Range("L112").Value = 210808.008
$ws.Range("M112").Value = -6962.0
$ws.Range("N112").Value = -213024.008
# Row 132
$ws.Range("H132").Value = 2756.5908
$ws.Range("I132").Value = 2745.0
$ws.Range("J132").Value = 3000.0
$ws.Range("K132").Value = 8235.0
$ws.Range("L132").Value = 9000.0
$ws.Range("M132").Value = -5705.0
$ws.Range("N132").Value = -14060.0
# Row 137
$ws.Range("H137").Value = 2113.9348
$ws.Range("I137").Value = 1521.8125
$ws.Range("J137").Value = 3467.3572
$ws.Range("K137").Value = 4565.4375
$ws.Range("L137").Value = 10402.0716
$ws.Range("M137").Value = -2015.4375
$ws.Range("N137").Value = -15502.0716
# Row 138
$ws.Range("H138").Value = 3646.861
$ws.Range("I138").Value = 3537.3572
$ws.Range("J138").Value = 3673.2932
$ws.Range("K138").Value = 10612.0716
$ws.Range("L138").Value = 11019.8796
$ws.Range("M138").Value = -5472.071599999999
$ws.Range("N138").Value = -21299.8796
# Row 141
$ws.Range("H141").Value = 2663.1428
$ws.Range("I141").Value = 2557.0557
$ws.Range("J141").Value = 3299.6667
$ws.Range("K141").Value = 7671.1671
$ws.Range("L141").Value = 9899.000100000001
$ws.Range("M141").Value = -2491.1671
$ws.Range("N141").Value = -20259.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 500.5
$ws.Range("I5").Value = 0.0
$ws.Range("J5").Value = 500.5
$ws.Range("K5").Value = 0.0
$ws.Range("L5").Value = 500.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -724.5
# Row 32
$ws.Range("H32").Value = 8558.06
$ws.Range("I32").Value = 5256.7744
$ws.Range("J32").Value = 16640.518
$ws.Range("K32").Value = 5256.7744
$ws.Range("L32").Value = 16640.518
$ws.Range("M32").Value = -4969.7744
$ws.Range("N32").Value = -17214.518
# Row 34
$ws.Range("H34").Value = 0.0
$ws.Range("I34").Value = 0.0
$ws.Range("K34").Value = 0.0
$ws.Range("M34").ClearContents()
# Row 45
$ws.Range("H45").Value = 1198.75
$ws.Range("I45").Value = 953.0
$ws.Range("K45").Value = 953.0
$ws.Range("M45").Value = -576.0
# Row 61
$ws.Range("H61").Value = 68969336.0
$ws.Range("I61").Value = 71432150.0
$ws.Range("J61").Value = 10500.0
$ws.Range("K61").Value = 71432150.0
$ws.Range("L61").Value = 10500.0
$ws.Range("M61").Value = -71431938.0
$ws.Range("N61").Value = -10924.0
# Row 74
$ws.Range("H74").Value = 37042584.0
$ws.Range("I74").Value = 41672604.0
$ws.Range("J74").Value = 2416.6667
$ws.Range("K74").Value = 41672604.0
$ws.Range("L74").Value = 2416.6667
$ws.Range("M74").Value = -41671730.0
$ws.Range("N74").Value = -4164.6667
# Row 77
$ws.Range("H77").Value = 37042584.0
$ws.Range("I77").Value = 41672604.0
$ws.Range("J77").Value = 2416.6667
$ws.Range("K77").Value = 208363020.0
$ws.Range("L77").Value = 12083.3335
$ws.Range("M77").Value = -208358652.0
$ws.Range("N77").Value = -20819.3335
# Row 96
$ws.Range("H96").Value = 25332.666
$ws.Range("J96").Value = 25332.666
$ws.Range("L96").Value = 25332.666
$ws.Range("N96").Value = -30824.666
# Row 97
$ws.Range("H97").Value = 605.7778
$ws.Range("I97").Value = 632.7059
$ws.Range("J97").Value = 148.0
$ws.Range("K97").Value = 632.7059
$ws.Range("L97").Value = 148.0
$ws.Range("M97").Value = -136.7059
$ws.Range("N97").Value = -1140.0
# Row 102
$ws.Range("H102").Value = 4631384.0
$ws.Range("I102").Value = 5292331.5
$ws.Range("K102").Value = 5292331.5
$ws.Range("M102").Value = -5290709.5
# Row 110
$ws.Range("H110").Value = 74491.86
$ws.Range("I110").Value = 113619.445
$ws.Range("K110").Value = 113619.445
$ws.Range("M110").Value = -111574.445
# Row 126
$ws.Range("H126").Value = 8000.0
$ws.Range("I126").Value = 8000.0
$ws.Range("K126").Value = 24000.0
$ws.Range("M126").Value = -21530.0
# Row 132
$ws.Range("H132").Value = 3850915.5
$ws.Range("I132").Value = 4004552.0
$ws.Range("J132").Value = 10000.0
$ws.Range("K132").Value = 12013656.0
$ws.Range("L132").Value = 30000.0
$ws.Range("M132").Value = -12011126.0
$ws.Range("N132").Value = -35060.0
# Row 136
$ws.Range("H136").Value = 68969336.0
$ws.Range("I136").Value = 71432150.0
$ws.Range("J136").Value = 10500.0
$ws.Range("K136").Value = 214296450.0
$ws.Range("L136").Value = 31500.0
$ws.Range("M136").Value = -214293900.0
$ws.Range("N136").Value = -36600.0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 500.5
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = 500.5
$ws.Range("K4").Value = 0.0
$ws.Range("L4").Value = 500.5
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -730.5
# Row 20
$ws.Range("H20").Value = 1791.1923
$ws.Range("I20").Value = 1677.9474
$ws.Range("J20").Value = 2098.5715
$ws.Range("K20").Value = 1677.9474
$ws.Range("L20").Value = 2098.5715
$ws.Range("M20").Value = -1430.9474
$ws.Range("N20").Value = -2592.5715
# Row 75
$ws.Range("H75").Value = 29166.666
$ws.Range("I75").Value = 7500.0
$ws.Range("K75").Value = 7500.0
$ws.Range("M75").Value = -6564.0
# Row 76
$ws.Range("H76").Value = 0.0
$ws.Range("J76").Value = 0.0
$ws.Range("L76").Value = 0.0
$ws.Range("N76").ClearContents()
# Row 78
$ws.Range("H78").Value = 29166.666
$ws.Range("I78").Value = 7500.0
$ws.Range("K78").Value = 22500.0
$ws.Range("M78").Value = -17820.0
# Row 79
$ws.Range("H79").Value = 0.0
$ws.Range("J79").Value = 0.0
$ws.Range("L79").Value = 0.0
$ws.Range("N79").ClearContents()
# Row 80
$ws.Range("H80").Value = 50001092.0
$ws.Range("I80").Value = 1175.6
$ws.Range("J80").Value = 100001010.0
$ws.Range("K80").Value = 1175.6
$ws.Range("L80").Value = 100001010.0
$ws.Range("M80").Value = -177.5999999999999
$ws.Range("N80").Value = -100003006.0
# Row 82
$ws.Range("H82").Value = 11598.4
$ws.Range("I82").Value = 11598.4
$ws.Range("K82").Value = 11598.4
$ws.Range("M82").Value = -11215.4
# Row 83
$ws.Range("H83").Value = 50001092.0
$ws.Range("I83").Value = 1175.6
$ws.Range("J83").Value = 100001010.0
$ws.Range("K83").Value = 5878.0
$ws.Range("L83").Value = 500005050.0
$ws.Range("M83").Value = -886.0
$ws.Range("N83").Value = -500015034.0
# Row 85
$ws.Range("H85").Value = 11598.4
$ws.Range("I85").Value = 11598.4
$ws.Range("K85").Value = 11598.4
$ws.Range("M85").Value = -10272.4
# Row 86
$ws.Range("H86").Value = 11624.75
$ws.Range("I86").Value = 13750.0
$ws.Range("J86").Value = 9499.5
$ws.Range("K86").Value = 13750.0
$ws.Range("L86").Value = 9499.5
$ws.Range("M86").Value = -12627.0
$ws.Range("N86").Value = -11745.5
# Row 88
$ws.Range("H88").Value = 29999.0
$ws.Range("J88").Value = 29999.0
$ws.Range("L88").Value = 29999.0
$ws.Range("N88").Value = -30811.0
# Row 89
$ws.Range("H89").Value = 11624.75
$ws.Range("I89").Value = 13750.0
$ws.Range("J89").Value = 9499.5
$ws.Range("K89").Value = 68750.0
$ws.Range("L89").Value = 47497.5
$ws.Range("M89").Value = -63134.0
$ws.Range("N89").Value = -58729.5
# Row 91
$ws.Range("H91").Value = 29999.0
$ws.Range("J91").Value = 29999.0
$ws.Range("L91").Value = 29999.0
$ws.Range("N91").Value = -32807.0
# Row 94
$ws.Range("H94").Value = 747.0
$ws.Range("I94").Value = 731.1875
$ws.Range("K94").Value = 731.1875
$ws.Range("M94").Value = -280.1875
# Row 134
$ws.Range("H134").Value = 13516420.0
$ws.Range("I134").Value = 14288415.0
$ws.Range("J134").Value = 6499.5
$ws.Range("K134").Value = 42865245.0
$ws.Range("L134").Value = 19498.5
$ws.Range("M134").Value = -42862710.0
$ws.Range("N134").Value = -24568.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 11871.444
$ws.Range("I22").Value = 14834.857
$ws.Range("J22").Value = 1499.5
$ws.Range("K22").Value = 14834.857
$ws.Range("L22").Value = 1499.5
$ws.Range("M22").Value = -14484.857
$ws.Range("N22").Value = -2199.5
# Row 31
$ws.Range("H31").Value = 5841.0
$ws.Range("I31").Value = 4601.1055
$ws.Range("J31").Value = 7653.154
$ws.Range("K31").Value = 4601.1055
$ws.Range("L31").Value = 7653.154
$ws.Range("M31").Value = -4306.1055
$ws.Range("N31").Value = -8243.154
# Row 34
$ws.Range("H34").Value = 5841.0
$ws.Range("I34").Value = 4601.1055
$ws.Range("J34").Value = 7653.154
$ws.Range("K34").Value = 4601.1055
$ws.Range("L34").Value = 7653.154
$ws.Range("M34").Value = -4399.1055
$ws.Range("N34").Value = -8057.154
# Row 43
$ws.Range("H43").Value = 23799.4
$ws.Range("J43").Value = 27499.25
$ws.Range("L43").Value = 27499.25
$ws.Range("N43").Value = -27867.25
# Row 50
$ws.Range("H50").Value = 0.0
$ws.Range("I50").Value = 0.0
$ws.Range("K50").Value = 0.0
$ws.Range("M50").ClearContents()
# Row 58
$ws.Range("H58").Value = 29419278.0
$ws.Range("I58").Value = 29419278.0
$ws.Range("K58").Value = 29419278.0
$ws.Range("M58").Value = -29419075.0
# Row 62
$ws.Range("H62").Value = 7479.0
$ws.Range("I62").Value = 7479.0
$ws.Range("J62").Value = 0.0
$ws.Range("K62").Value = 7479.0
$ws.Range("L62").Value = 0.0
$ws.Range("M62").Value = -6855.0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 7479.0
$ws.Range("I65").Value = 7479.0
$ws.Range("J65").Value = 0.0
$ws.Range("K65").Value = 37395.0
$ws.Range("L65").Value = 0.0
$ws.Range("M65").Value = -34275.0
$ws.Range("N65").ClearContents()
# Row 101
$ws.Range("H101").Value = 23799.4
$ws.Range("J101").Value = 27499.25
$ws.Range("L101").Value = 27499.25
$ws.Range("N101").Value = -33989.25
# Row 107
$ws.Range("H107").Value = 802446.5
$ws.Range("I107").Value = 1134141.5
$ws.Range("J107").Value = 205395.6
$ws.Range("K107").Value = 1134141.5
$ws.Range("L107").Value = 205395.6
$ws.Range("M107").Value = -1132221.5
$ws.Range("N107").Value = -209235.6
# Row 127
$ws.Range("H127").Value = 104999.5
$ws.Range("J127").Value = 104999.5
$ws.Range("L127").Value = 104999.5
$ws.Range("N127").Value = -114919.5
# Row 132
$ws.Range("H132").Value = 16668334.0
$ws.Range("I132").Value = 18520208.0
$ws.Range("J132").Value = 1465.5
$ws.Range("K132").Value = 55560624.0
$ws.Range("L132").Value = 4396.5
$ws.Range("M132").Value = -55558094.0
$ws.Range("N132").Value = -9456.5
# Row 136
$ws.Range("H136").Value = 29419278.0
$ws.Range("I136").Value = 29419278.0
$ws.Range("K136").Value = 88257834.0
$ws.Range("M136").Value = -88255284.0
# Row 141
$ws.Range("H141").Value = 55157.867
$ws.Range("J141").Value = 57963.11
$ws.Range("L141").Value = 57963.11
$ws.Range("N141").Value = -68323.11

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 239.26315
$ws.Range("J12").Value = 209.15384
$ws.Range("L12").Value = 627.4615200000001
$ws.Range("N12").Value = -973.4615200000001
# Row 14
$ws.Range("H14").Value = 515.6667
$ws.Range("I14").Value = 515.6667
$ws.Range("K14").Value = 1547.0001
$ws.Range("M14").Value = -1374.0001
# Row 52
$ws.Range("H52").Value = 2819.6
$ws.Range("J52").Value = 2819.6
$ws.Range("L52").Value = 8458.8
$ws.Range("N52").Value = -8990.8
# Row 60
$ws.Range("H60").Value = 5548.4614
$ws.Range("I60").Value = 326.0
$ws.Range("K60").Value = 978.0
$ws.Range("M60").Value = -727.0
# Row 94
$ws.Range("H94").Value = 16502.75
$ws.Range("J94").Value = 18574.572
$ws.Range("L94").Value = 55723.716
$ws.Range("N94").Value = -57075.716
# Row 98
$ws.Range("H98").Value = 1018.5
$ws.Range("I98").Value = 1061.0834
$ws.Range("J98").Value = 933.3333
$ws.Range("K98").Value = 3183.2502
$ws.Range("L98").Value = 2799.9999
$ws.Range("M98").Value = -1685.2502
$ws.Range("N98").Value = -5795.9999
# Row 122
$ws.Range("H122").Value = 1805.9286
$ws.Range("I122").Value = 1030.1818
$ws.Range("K122").Value = 9271.6362
$ws.Range("M122").Value = -6821.636200000001
# Row 132
$ws.Range("H132").Value = 2196.0
$ws.Range("I132").Value = 2299.1428
$ws.Range("J132").Value = 2105.75
$ws.Range("K132").Value = 20692.2852
$ws.Range("L132").Value = 18951.75
$ws.Range("M132").Value = -18162.2852
$ws.Range("N132").Value = -24011.75
# Row 140
$ws.Range("H140").Value = 2872.4
$ws.Range("I140").Value = 2872.4
$ws.Range("K140").Value = 8617.2
$ws.Range("M140").Value = -3437.200000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 147.16667
$ws.Range("J2").Value = 275.0
$ws.Range("L2").Value = 275.0
$ws.Range("N2").Value = -501.0
# Row 64
$ws.Range("H64").Value = 70000.0
$ws.Range("J64").Value = 0.0
$ws.Range("L64").Value = 0.0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 70000.0
$ws.Range("J67").Value = 0.0
$ws.Range("L67").Value = 0.0
$ws.Range("N67").ClearContents()
# Row 80
$ws.Range("H80").Value = 5785.143
$ws.Range("I80").Value = 4374.25
$ws.Range("J80").Value = 7666.3335
$ws.Range("K80").Value = 4374.25
$ws.Range("L80").Value = 7666.3335
$ws.Range("M80").Value = -3376.25
$ws.Range("N80").Value = -9662.3335
# Row 83
$ws.Range("H83").Value = 5785.143
$ws.Range("I83").Value = 4374.25
$ws.Range("J83").Value = 7666.3335
$ws.Range("K83").Value = 21871.25
$ws.Range("L83").Value = 38331.6675
$ws.Range("M83").Value = -16879.25
$ws.Range("N83").Value = -48315.6675
# Row 97
$ws.Range("H97").Value = 1115.5
$ws.Range("I97").Value = 803.25
$ws.Range("J97").Value = 1427.75
$ws.Range("K97").Value = 803.25
$ws.Range("L97").Value = 1427.75
$ws.Range("M97").Value = -307.25
$ws.Range("N97").Value = -2419.75
# Row 102
$ws.Range("H102").Value = 3344.6
$ws.Range("I102").Value = 3180.875
$ws.Range("K102").Value = 3180.875
$ws.Range("M102").Value = -1558.875
# Row 113
$ws.Range("H113").Value = 132062.38
$ws.Range("I113").Value = 205499.6
$ws.Range("J113").Value = 9667.0
$ws.Range("K113").Value = 205499.6
$ws.Range("L113").Value = 9667.0
$ws.Range("M113").Value = -203329.6
$ws.Range("N113").Value = -14007.0
# Row 122
$ws.Range("H122").Value = 53228.582
$ws.Range("I122").Value = 72981.94
$ws.Range("J122").Value = 5256.143
$ws.Range("K122").Value = 218945.82
$ws.Range("L122").Value = 15768.429
$ws.Range("M122").Value = -216495.82
$ws.Range("N122").Value = -20668.429
# Row 123
$ws.Range("H123").Value = 67492.664
$ws.Range("J123").Value = 67492.664
$ws.Range("L123").Value = 67492.664
$ws.Range("N123").Value = -72392.664
# Row 125
$ws.Range("H125").Value = 46862.4
$ws.Range("J125").Value = 46328.0
$ws.Range("L125").Value = 46328.0
$ws.Range("N125").Value = -51248.0
# Row 126
$ws.Range("H126").Value = 7174.75
$ws.Range("I126").Value = 6763.9165
$ws.Range("J126").Value = 8407.25
$ws.Range("K126").Value = 20291.7495
$ws.Range("L126").Value = 25221.75
$ws.Range("M126").Value = -17821.7495
$ws.Range("N126").Value = -30161.75
# Row 132
$ws.Range("H132").Value = 3476532.2
$ws.Range("I132").Value = 3575612.0
$ws.Range("K132").Value = 10726836.0
$ws.Range("M132").Value = -10724306.0

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2709.577
$ws.Range("I16").Value = 1179.3334
$ws.Range("J16").Value = 4021.2144
$ws.Range("K16").Value = 1179.3334
$ws.Range("L16").Value = 4021.2144
$ws.Range("M16").Value = -1009.3334
$ws.Range("N16").Value = -4361.2144
# Row 22
$ws.Range("H22").Value = 1754.75
$ws.Range("I22").Value = 2285.2856
$ws.Range("K22").Value = 2285.2856
$ws.Range("M22").Value = -1990.2856
# Row 27
$ws.Range("H27").Value = 1754.75
$ws.Range("I27").Value = 2285.2856
$ws.Range("K27").Value = 2285.2856
$ws.Range("M27").Value = -2178.2856
# Row 30
$ws.Range("H30").Value = 3025.0
$ws.Range("I30").Value = 1112.5
$ws.Range("J30").Value = 14500.0
$ws.Range("K30").Value = 1112.5
$ws.Range("L30").Value = 14500.0
$ws.Range("M30").Value = -1004.5
$ws.Range("N30").Value = -14716.0
# Row 40
$ws.Range("H40").Value = 4498.6665
$ws.Range("J40").Value = 4998.6665
$ws.Range("L40").Value = 4998.6665
$ws.Range("N40").Value = -5270.6665
# Row 74
$ws.Range("H74").Value = 18849.5
$ws.Range("I74").Value = 18849.5
$ws.Range("K74").Value = 18849.5
$ws.Range("M74").Value = -17851.5
# Row 77
$ws.Range("H77").Value = 18849.5
$ws.Range("I77").Value = 18849.5
$ws.Range("K77").Value = 56548.5
$ws.Range("M77").Value = -51556.5
# Row 82
$ws.Range("H82").Value = 1456.5834
$ws.Range("I82").Value = 1397.7778
$ws.Range("J82").Value = 1633.0
$ws.Range("K82").Value = 1397.7778
$ws.Range("L82").Value = 1633.0
$ws.Range("M82").Value = -1036.7778
$ws.Range("N82").Value = -2355.0
# Row 85
$ws.Range("H85").Value = 1456.5834
$ws.Range("I85").Value = 1397.7778
$ws.Range("J85").Value = 1633.0
$ws.Range("K85").Value = 1397.7778
$ws.Range("L85").Value = 1633.0
$ws.Range("M85").Value = -149.7778000000001
$ws.Range("N85").Value = -4129.0
# Row 104
$ws.Range("H104").Value = 58554.5
$ws.Range("J104").Value = 58554.5
$ws.Range("L104").Value = 58554.5
$ws.Range("N104").Value = -65542.5
# Row 122
$ws.Range("H122").Value = 4253.857
$ws.Range("I122").Value = 4253.857
$ws.Range("K122").Value = 12761.571
$ws.Range("M122").Value = -10311.571
# Row 136
$ws.Range("H136").Value = 2858.25
$ws.Range("J136").Value = 3004.5
$ws.Range("L136").Value = 9013.5
$ws.Range("N136").Value = -14113.5
# Row 140
$ws.Range("H140").Value = 95265.0
$ws.Range("J140").Value = 95265.0
$ws.Range("L140").Value = 95265.0
$ws.Range("N140").Value = -105625.0

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1580.9445
$ws.Range("I100").Value = 1520.4667
$ws.Range("K100").Value = 3040.9334
$ws.Range("M100").Value = -2499.9334
# Row 107
$ws.Range("H107").Value = 1766.5
$ws.Range("I107").Value = 444.66666
$ws.Range("K107").Value = 1333.99998
$ws.Range("M107").Value = 586.00002
# Row 122
$ws.Range("H122").Value = 3325.15
$ws.Range("I122").Value = 3076.4119
$ws.Range("J122").Value = 4734.6665
$ws.Range("K122").Value = 9229.235700000001
$ws.Range("L122").Value = 14203.9995
$ws.Range("M122").Value = -6779.235700000001
$ws.Range("N122").Value = -19103.9995
# Row 126
$ws.Range("H126").Value = 2665.7144
$ws.Range("I126").Value = 2665.7144
$ws.Range("K126").Value = 7997.1432
$ws.Range("M126").Value = -5527.1432
# Row 132
$ws.Range("H132").Value = 16135983.0
$ws.Range("I132").Value = 21743344.0
$ws.Range("J132").Value = 14823.125
$ws.Range("K132").Value = 65230032.0
$ws.Range("L132").Value = 44469.375
$ws.Range("M132").Value = -65227502.0
$ws.Range("N132").Value = -49529.375
# Row 135
$ws.Range("H135").Value = 119805.89
$ws.Range("J135").Value = 119805.89
$ws.Range("L135").Value = 119805.89
$ws.Range("N135").Value = -129945.89
# Row 136
$ws.Range("H136").Value = 13514827.0
$ws.Range("I136").Value = 14286817.0
$ws.Range("K136").Value = 42860451.0
$ws.Range("M136").Value = -42857901.0
# Row 140
$ws.Range("H140").Value = 69517.0
$ws.Range("J140").Value = 74892.664
$ws.Range("L140").Value = 74892.664
$ws.Range("N140").Value = -85252.664

